$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column E (rows 2-12) from 50 to 70
$ws.Range("E2:E12").Value = 70
